# Rename the second sheet from "Sheet2" to "Shortcut Keys"
$wb = $excel.ActiveWorkbook
$sheet2 = $wb.Sheets.Item("Sheet2")
$sheet2.Name = "Shortcut Keys"

# Fill in the lesson names for Episodes 13-17 (rows 14-18) on the
# "Completed Lessons" sheet
$ws = $wb.Sheets.Item("Completed Lessons")
$ws.Range("C14").Value = "Entering Text to Create Spreadsheet Titles"
$ws.Range("C15").Value = "Working with Numeric Data in Excel"
$ws.Range("C16").Value = "Entering Date Values in Excel"
$ws.Range("C17").Value = "Working with Cell References"
$ws.Range("C18").Value = "Creating Basic Formulas in Excel"

# Update the active selection to C19, matching the saved view state
$ws.Range("C19").Select()
